$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------------
# 1) Sheet3: bump the reference date and paste in the new day's raw lookup
#    table (B20:B36). This mirrors the daily workflow: new readings land in
#    the lookup block, the VLOOKUP table date advances by one day.
# ---------------------------------------------------------------------------
$ws3.Range("C1").Value = 45226

$newLookup = @{
    20 = 17.066972185230046
    21 = 1.4844142670901579
    22 = 12.072014613499624
    23 = 0
    24 = 0
    25 = 0
    26 = 4.8163953117382547
    27 = 8.2237724383598074
    28 = 5.7780022334884098
    29 = 5.4813282516161346
    30 = 1.4302902346381163
    31 = 4.654328693117928
    32 = 9.0140083746758908
    33 = 7.0412601452528989
    34 = 4.0032998851250161
    35 = 7.170416399557876
    36 = 41.535785349703858
}
foreach ($r in 20..36) {
    $ws3.Range("B$r").Value = $newLookup[$r]
}

# ---------------------------------------------------------------------------
# 2) Sheet1: make the CB VLOOKUPs anchor the lookup key column ($B) so the
#    formula can be filled right into the new CC column unchanged.
# ---------------------------------------------------------------------------
foreach ($r in 2..18) {
    $ws1.Range("CB$r").Formula = '=VLOOKUP($B' + $r + ',Sheet3!$B$1:$C$18,2,)'
}

# ---------------------------------------------------------------------------
# 3) Sheet3: freeze today's VLOOKUP results (C2:C18) to plain values now that
#    they have recalculated against the new lookup table above.
# ---------------------------------------------------------------------------
foreach ($r in 2..18) {
    $cell = $ws3.Range("C$r")
    $cell.Value = $cell.Value2
}

# ---------------------------------------------------------------------------
# 4) Sheet1: add the new day's column (CC) - header label plus the same
#    VLOOKUP carried across from CB.
# ---------------------------------------------------------------------------
$ws1.Range("CC1").Value = "27-oct"
foreach ($r in 2..18) {
    $ws1.Range("CC$r").Formula = '=VLOOKUP($B' + $r + ',Sheet3!$B$1:$C$18,2,)'
}

# ---------------------------------------------------------------------------
# 5) View state: Sheet3 keeps its own remembered scroll/selection (set while
#    it is active), then Sheet1 is (re)activated last so it ends up as the
#    selected tab with its own remembered selection.
# ---------------------------------------------------------------------------
$ws3.Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws3.Range("B33").Select()

$ws1.Select()
$ws1.Range("CF8").Select()
